# "Added after the seminar"
# Adds a Pass/Fail style "Pass" column value to the login sheet and
# populates two new rows (Paul/PL, Sam/US) on the info sheet.

$wb = $excel.ActiveWorkbook

# --- login sheet: D2 gets a new "Pass" label ---
$ws1 = $wb.Worksheets.Item("login")
$ws1.Range("D2").Value = "Pass"

# --- info sheet: add rows 2 and 3 ---
$ws2 = $wb.Worksheets.Item("info")
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Paul"
$ws2.Range("C2").Value = "PL"
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Sam"
$ws2.Range("C3").Value = "US"

# Mirror the author's final selection state: info sheet has A2:C3
# selected, but the login sheet remains the active tab.
[void]$ws2.Range("A2:C3").Select()
$ws1.Activate()
